# Rename the two existing sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "nhap-thanhpham"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "xuat-thanhpham"

# Add a third sheet at the end (clone the first sheet so the column /
# page-setup metadata matches what a real worksheet in this workbook looks
# like, then overwrite its contents below)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1.Copy($null, $lastSheet)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "ton-thanhpham"

# ---------------------------------------------------------------
# Sheet 1: nhap-thanhpham
# ---------------------------------------------------------------
$ws1.Range("A1").Value = "Tên Hàng"
$ws1.Range("B1").Value = "MCU"
$ws1.Range("C1").Value = "Sổ Hợp Đồng"
$ws1.Range("D1").Value = "Chip"
$ws1.Range("E1").Value = "Ngày Nhập"
$ws1.Range("F1").Value = "Số Lượng"

$ws1.Range("A2").Value = "smooth"
$ws1.Range("B2").Value = "mcu"
$ws1.Range("C2").Value = "da"
$ws1.Range("D2").Value = "chip"

$ws1.Range("E2").Formula = '="2021-10-05"'
$ws1.Range("E2").Copy()
$ws1.Range("E2").PasteSpecial(-4163)

$ws1.Range("F2").Value = 50

$ws1.Range("G1:J2").ClearContents()
$ws1.Columns("G:J").Delete()

# ---------------------------------------------------------------
# Sheet 2: xuat-thanhpham
# ---------------------------------------------------------------
$ws2.Range("A1").Value = "Tên Hàng"
$ws2.Range("B1").Value = "MCU"
$ws2.Range("C1").Value = "Sổ Hợp Đồng"
$ws2.Range("D1").Value = "Chip"
$ws2.Range("E1").Value = "Ngày Nhập"
$ws2.Range("F1").Value = "Số Lượng"

$ws2.Range("A2").Value = "smooth"
$ws2.Range("B2").Value = "mcu"
$ws2.Range("C2").Value = "da"
$ws2.Range("D2").Value = "chip"

$ws2.Range("E2").Formula = '="2021-10-06"'
$ws2.Range("E2").Copy()
$ws2.Range("E2").PasteSpecial(-4163)

$ws2.Range("F2").Value = 40

$ws2.Range("G1:J1").ClearContents()
$ws2.Columns("G:J").Delete()

# ---------------------------------------------------------------
# Sheet 3: ton-thanhpham
# ---------------------------------------------------------------
$ws3.Range("A1").Value = "Tên Hàng"
$ws3.Range("B1").Value = "Số Lượng"
$ws3.Range("C1").Value = "Đơn Vị Tính"

$ws3.Range("A2").Value = "smooth"
$ws3.Range("B2").Value = 10
$ws3.Range("C2").Value = "none"

$ws3.Range("D1:J2").ClearContents()
$ws3.Columns("D:J").Delete()

Write-Host "done"
